$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.480.12'
$ws.Range('E2').Value = '  -1.34%  '
$ws.Range('D3').Value = '2.045.48'
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '246.72'
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '54.08'
$ws.Range('E8').Value = '  -6.40%  '
$ws.Range('D9').Value = '62.06'
$ws.Range('E9').Value = '  +4.30%  '
$ws.Range('D10').Value = '0.362'
$ws.Range('E10').Value = '  -3.02%  '
$ws.Range('D11').Value = '0.0744'
$ws.Range('E11').Value = '  -4.76%  '
$ws.Range('E12').Value = '  -3.69%  '
$ws.Range('D13').Value = '0.942'
$ws.Range('E13').Value = '  +7.45%  '
$ws.Range('D14').Value = '14.60'
$ws.Range('E14').Value = '  -4.47%  '
$ws.Range('D15').Value = '2.348.22'
$ws.Range('E15').Value = '  -0.12%  '
$ws.Range('D16').Value = '5.40'
$ws.Range('E16').Value = '  -4.65%  '
$ws.Range('D17').Value = '2.045.40'
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').Value = '36.377.09'
$ws.Range('E18').Value = '  -1.53%  '
$ws.Range('D19').Value = '17.01'
$ws.Range('E19').Value = '  -5.16%  '
$ws.Range('D20').Value = '71.55'
$ws.Range('E20').Value = '  -2.76%  '
$ws.Range('D21').Value = '0.0₃0853'
$ws.Range('E21').Value = '  -4.59%  '
$ws.Range('D22').Value = '236.88'
$ws.Range('E22').Value = '  +0.44%  '
$ws.Range('D23').Value = '5.18'
$ws.Range('E23').Value = '  -4.57%  '
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('E25').Value = '  -2.69%  '
$ws.Range('E26').Value = '  +3.18%  '
$ws.Range('D27').Value = '164.74'
$ws.Range('E27').Value = '  -2.60%  '
$ws.Range('D28').Value = '9.14'
$ws.Range('E28').Value = '  -11.37%  '
$ws.Range('D29').Value = '19.81'
$ws.Range('E29').Value = '  -1.12%  '
$ws.Range('D31').Value = '1.18'
$ws.Range('E31').Value = '  +6.38%  '
$ws.Range('D32').Value = '5.02'
$ws.Range('E32').Value = '  -10.10%  '
$ws.Range('D33').Value = '4.41'
$ws.Range('E33').Value = '  -6.67%  '
$ws.Range('D34').Value = '0.0590'
$ws.Range('E34').Value = '  -4.36%  '
$ws.Range('D35').Value = '0.0878'
$ws.Range('E35').Value = '  +7.89%  '
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('E37').Value = '  -0.77%  '
$ws.Range('D38').Value = '2.20'
$ws.Range('E38').Value = '  -6.12%  '
$ws.Range('D39').Value = '5.02'
$ws.Range('E39').Value = '  -2.85%  '
$ws.Range('E40').Value = '  -6.65%  '
$ws.Range('E41').Value = '  -4.85%  '
$ws.Range('E42').Value = '  -4.89%  '
$ws.Range('E43').Value = '  -4.78%  '
$ws.Range('D44').Value = '93.59'
$ws.Range('E44').Value = '  -3.67%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '1.400.53'
$ws.Range('E45').Value = '  +7.15%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').Value = '0.0899'
$ws.Range('E46').Value = '  -5.93%  '
$ws.Range('D47').Value = '15.76'
$ws.Range('E47').Value = '  -6.85%  '
$ws.Range('D48').Value = '7.37'
$ws.Range('E48').Value = '  +9.27%  '
$ws.Range('E49').Value = '  +1.64%  '
$ws.Range('D50').Value = '2.25'
$ws.Range('E50').Value = '  -4.36%  '
$ws.Range('D51').Value = '2.230.42'
$ws.Range('E51').Value = '  -0.27%  '
